$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "as of" snapshot timestamp in column D for all data rows (2-57)
$ws.Range("D2:D57").Value = 45966.40896990741

# Update rows 18-51: charging station (A), terminal name (B), last charge end time (C)
$ws.Range("A18").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B18").Value = "101号直流"
$ws.Range("C18").Value = 45954.028229166666
$ws.Range("A19").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B19").Value = "008B号直流"
$ws.Range("C19").Value = 45959.55945601852
$ws.Range("A20").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B20").Value = "904号直流"
$ws.Range("C20").Value = 45962.643703703703
$ws.Range("A21").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B21").Value = "306号直流"
$ws.Range("C21").Value = 45964.263055555559
$ws.Range("A22").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B22").Value = "004A号直流"
$ws.Range("C22").Value = 45964.528668981482
$ws.Range("A23").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B23").Value = "505号直流"
$ws.Range("C23").Value = 45964.534479166665
$ws.Range("A24").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B24").Value = "A02号直流"
$ws.Range("C24").Value = 45965.056006944447
$ws.Range("A25").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B25").Value = "502号直流"
$ws.Range("C25").Value = 45965.254895833335
$ws.Range("A26").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B26").Value = "702号直流"
$ws.Range("C26").Value = 45965.261817129627
$ws.Range("A27").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B27").Value = "402号直流"
$ws.Range("C27").Value = 45965.27380787037
$ws.Range("A28").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B28").Value = "703号直流"
$ws.Range("C28").Value = 45965.306504629632
$ws.Range("A29").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B29").Value = "203号直流"
$ws.Range("C29").Value = 45965.339097222219
$ws.Range("A30").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B30").Value = "A01号直流"
$ws.Range("C30").Value = 45965.373263888891
$ws.Range("A31").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B31").Value = "B02号直流"
$ws.Range("C31").Value = 45965.519780092596
$ws.Range("A32").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B32").Value = "502号直流"
$ws.Range("C32").Value = 45965.528715277775
$ws.Range("A33").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B33").Value = "001B号直流"
$ws.Range("C33").Value = 45965.533333333333
$ws.Range("A34").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B34").Value = "503号直流"
$ws.Range("C34").Value = 45965.534953703704
$ws.Range("A35").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B35").Value = "405号直流"
$ws.Range("C35").Value = 45965.539178240739
$ws.Range("A36").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B36").Value = "805号直流"
$ws.Range("C36").Value = 45965.550335648149
$ws.Range("A37").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B37").Value = "103号直流"
$ws.Range("C37").Value = 45965.551053240742
$ws.Range("A38").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B38").Value = "602号直流"
$ws.Range("C38").Value = 45965.561493055553
$ws.Range("A39").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B39").Value = "804号直流"
$ws.Range("C39").Value = 45965.564560185187
$ws.Range("A40").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B40").Value = "201号直流"
$ws.Range("C40").Value = 45965.565578703703
$ws.Range("A41").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B41").Value = "604号直流"
$ws.Range("C41").Value = 45965.565891203703
$ws.Range("A42").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B42").Value = "802号直流"
$ws.Range("C42").Value = 45965.567743055559
$ws.Range("A43").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B43").Value = "009B号直流"
$ws.Range("C43").Value = 45965.582557870373
$ws.Range("A44").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B44").Value = "112号直流"
$ws.Range("C44").Value = 45965.587118055555
$ws.Range("A45").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B45").Value = "406号直流"
$ws.Range("C45").Value = 45965.587233796294
$ws.Range("A46").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B46").Value = "108号直流"
$ws.Range("C46").Value = 45965.620046296295
$ws.Range("A47").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B47").Value = "101号直流"
$ws.Range("C47").Value = 45965.622199074074
$ws.Range("A48").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B48").Value = "903号直流"
$ws.Range("C48").Value = 45965.638831018521
$ws.Range("A49").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B49").Value = "201号直流"
$ws.Range("C49").Value = 45965.661550925928
$ws.Range("A50").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B50").Value = "104号直流"
$ws.Range("C50").Value = 45965.679039351853
$ws.Range("A51").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B51").Value = "301号直流"
$ws.Range("C51").Value = 45965.74291666667

# Rows 52-57 no longer have data; clear their contents but keep formatting
$ws.Range("A52:D57").ClearContents()

# Update the active cell selection
$ws.Range("E16").Select()
